$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 116, pushing the existing row 116 (and below) down.
$ws.Rows.Item(116).Insert()

# Copy the number format (style) used by the date column from the row below
# (now row 117, the old row 116) onto the new row's date cell so it keeps the
# same date formatting.
$ws.Cells.Item(117, 4).Copy()
$ws.Cells.Item(116, 4).PasteSpecial(-4122)  # xlPasteFormats

# Populate the newly inserted row 116 with the new record's values.
$ws.Cells.Item(116, 1).Value = 11
$ws.Cells.Item(116, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(116, 3).Value = "Bíobío"
$ws.Cells.Item(116, 4).Value = 44656
$ws.Cells.Item(116, 5).Value = 8
$ws.Cells.Item(116, 6).Value = 100112043
$ws.Cells.Item(116, 7).Value = "Pepino ensalada"
$ws.Cells.Item(116, 8).Value = "Sin especificar"
$ws.Cells.Item(116, 9).Value = "Primera"
$ws.Cells.Item(116, 10).Value = 240
$ws.Cells.Item(116, 11).Value = 12000
$ws.Cells.Item(116, 12).Value = 13000
$ws.Cells.Item(116, 13).Value = 12500
$ws.Cells.Item(116, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(116, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(116, 16).Value = 208
$ws.Cells.Item(116, 17).Value = 60
$ws.Cells.Item(116, 18).Value = "Hortaliza"
